$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("L5").Value = 0
$ws1.Range("M5").Value = 0
$ws1.Range("M6").Value = 0
$ws1.Range("I11").Value = 0
$ws1.Range("L12").Value = 0
$ws1.Range("L22").Value = 0
$ws1.Range("I25").Value = 0
$ws1.Range("D26").Value = 0
$ws1.Range("D28").Value = 0
$ws1.Range("L28").Value = 0
$ws1.Range("M28").Value = 0
$ws1.Range("D29").Value = 0
$ws1.Range("H29").Value = 0
$ws1.Range("I29").Value = 0
$ws1.Range("L29").Value = 0
$ws1.Range("M29").Value = 0
$ws1.Range("I30").Value = 0
$ws1.Range("M30").Value = 0
$ws1.Range("M31").Value = 0
$ws1.Range("C37").Value = 0
$ws1.Range("L37").Value = 0
$ws1.Range("M37").Value = 0
$ws1.Range("M38").Value = 0
$ws1.Range("K40").Value = 0
$ws1.Range("Q40").Value = 0
$ws1.Range("Q42").Value = 0
$ws1.Range("L43").Value = 0
$ws1.Range("M43").Value = 0
$ws1.Range("D44").Value = 0
$ws1.Range("I44").Value = 0
$ws1.Range("M44").Value = 0
$ws1.Range("N44").Value = 0
$ws1.Range("D45").Value = 0
$ws1.Range("L45").Value = 0
$ws1.Range("M45").Value = 0
$ws1.Range("D46").Value = 0
$ws1.Range("H46").Value = 0
$ws1.Range("I46").Value = 0
$ws1.Range("L46").Value = 0
$ws1.Range("M46").Value = 0
$ws1.Range("M48").Value = 0
$ws1.Range("L50").Value = 0
$ws1.Range("G51").Value = 0
$ws1.Range("M51").Value = 0
$ws1.Range("N51").Value = 0
$ws1.Range("Q54").Value = 0
$ws1.Range("C56").Value = '0 de 54'
$ws1.Range("D56").Value = '0 de 54'
$ws1.Range("G56").Value = '0 de 54'
$ws1.Range("H56").Value = '0 de 54'
$ws1.Range("I56").Value = '0 de 54'
$ws1.Range("K56").Value = '0 de 54'
$ws1.Range("L56").Value = '0 de 54'
$ws1.Range("M56").Value = '0 de 54'
$ws1.Range("N56").Value = '0 de 54'
$ws1.Range("Q56").Value = '0 de 54'

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("C1").Value = 'mayo'
$ws2.Range("D1").Value = 'junio'
$ws2.Range("E1").Value = 'julio'
$ws2.Range("F1").Value = 'agosto'
$ws2.Range("C5").Value = 4158.27
$ws2.Range("D5").Value = 5087.23
$ws2.Range("E5").Value = 9991.16
$ws2.Range("F5").Value = 0
$ws2.Range("C6").Value = 1528.39
$ws2.Range("D6").Value = 1516.28
$ws2.Range("E6").Value = 1795.71
$ws2.Range("F6").Value = 0
$ws2.Range("C11").Value = 0
$ws2.Range("E11").Value = 890.8
$ws2.Range("F11").Value = 0
$ws2.Range("D12").Value = 135.9
$ws2.Range("E12").Value = 320.98
$ws2.Range("F12").Value = 0
$ws2.Range("C14").Value = 2169.75
$ws2.Range("D14").Value = 456.84
$ws2.Range("E14").Value = 0
$ws2.Range("C19").Value = 0
$ws2.Range("D20").Value = 259.49
$ws2.Range("E20").Value = 0
$ws2.Range("C22").Value = 4141.42
$ws2.Range("D22").Value = 2733.67
$ws2.Range("E22").Value = 1710.72
$ws2.Range("F22").Value = 0
$ws2.Range("C25").Value = 61.78
$ws2.Range("D25").Value = 3114.83
$ws2.Range("E25").Value = 818.0599999999999
$ws2.Range("F25").Value = 0
$ws2.Range("C26").Value = 0
$ws2.Range("E26").Value = 1373.76
$ws2.Range("F26").Value = 0
$ws2.Range("D27").Value = 727.97
$ws2.Range("E27").Value = 0
$ws2.Range("C28").Value = 6249.76
$ws2.Range("D28").Value = 36680.13
$ws2.Range("E28").Value = 17469.82
$ws2.Range("F28").Value = 0
$ws2.Range("C29").Value = 7315.29
$ws2.Range("D29").Value = 3247.24
$ws2.Range("E29").Value = 14529.6
$ws2.Range("F29").Value = 0
$ws2.Range("C30").Value = 3563.29
$ws2.Range("D30").Value = 1079.23
$ws2.Range("E30").Value = 5996.2
$ws2.Range("F30").Value = 0
$ws2.Range("E31").Value = 739.1
$ws2.Range("F31").Value = 0
$ws2.Range("C32").Value = 0
$ws2.Range("D32").Value = 2568.3
$ws2.Range("E32").Value = 0
$ws2.Range("C37").Value = 5704.92
$ws2.Range("D37").Value = 14177.18
$ws2.Range("E37").Value = 4677
$ws2.Range("F37").Value = 0
$ws2.Range("E38").Value = 1186.08
$ws2.Range("F38").Value = 0
$ws2.Range("C40").Value = 2403.41
$ws2.Range("D40").Value = 2172.6
$ws2.Range("E40").Value = 2874.67
$ws2.Range("F40").Value = 0
$ws2.Range("C42").Value = 2689.09
$ws2.Range("D42").Value = -39.67
$ws2.Range("E42").Value = 582.66
$ws2.Range("F42").Value = 0
$ws2.Range("C43").Value = 0
$ws2.Range("D43").Value = 86.5
$ws2.Range("E43").Value = 3269.52
$ws2.Range("F43").Value = 0
$ws2.Range("C44").Value = 7574.03
$ws2.Range("D44").Value = 6905.1
$ws2.Range("E44").Value = 5003.99
$ws2.Range("F44").Value = 0
$ws2.Range("C45").Value = 731.63
$ws2.Range("D45").Value = 3152.12
$ws2.Range("E45").Value = 6207.1
$ws2.Range("F45").Value = 0
$ws2.Range("C46").Value = 722.54
$ws2.Range("D46").Value = 158.83
$ws2.Range("E46").Value = 1712.88
$ws2.Range("F46").Value = 0
$ws2.Range("C48").Value = 798
$ws2.Range("D48").Value = 0
$ws2.Range("E48").Value = 154.28
$ws2.Range("F48").Value = 0
$ws2.Range("E50").Value = 380.16
$ws2.Range("F50").Value = 0
$ws2.Range("C51").Value = 4953.13
$ws2.Range("D51").Value = 5333.85
$ws2.Range("E51").Value = 1439.92
$ws2.Range("F51").Value = 0
$ws2.Range("C52").Value = -11.75
$ws2.Range("D52").Value = 0
$ws2.Range("E54").Value = 581.26
$ws2.Range("F54").Value = 0
$ws2.Range("C55").Value = 144
$ws2.Range("D55").Value = 0
$ws2.Range("C56").Value = 54896.95
$ws2.Range("D56").Value = 89553.62
$ws2.Range("E56").Value = 83705.42999999999
$ws2.Range("F56").Value = 0

# --- Column width adjustments on VENTA MENSUAL ---
# target integer widths: C=13, D=14, E=14(unchanged), F=12, G=17(unchanged)
$ws2.Columns.Item(3).ColumnWidth = 12.166666666666666  # -> stored width 13
$ws2.Columns.Item(4).ColumnWidth = 13.166666666666666  # -> stored width 14
$ws2.Columns.Item(6).ColumnWidth = 11.166666666666666  # -> stored width 12
